# ---- Auto-generated edit script ----
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 418
$ws1.Range("F3").Value = 573
$ws1.Range("F4").Value = 878
$ws1.Range("F5").Value = 614
$ws1.Range("F6").Value = 790
$ws1.Range("F7").Value = 360
$ws1.Range("F9").Value = 109
$ws1.Range("F10").Value = 1121
$ws1.Range("F11").Value = 570
$ws1.Range("F12").Value = 333
$ws1.Range("F13").Value = 448
$ws1.Range("F15").Value = 289

# Insert new row 19 (黑塔利亚Only) in sheet1, shifting old rows 19-21 down to 20-22
$ws1.Rows("19:19").Insert()
$ws1.Range("A20").Copy()
$ws1.Range("A19").PasteSpecial(-4122)
$ws1.Range("B20").Copy()
$ws1.Range("B19").PasteSpecial(-4122)
$ws1.Range("C20:I20").Copy()
$ws1.Range("C19").PasteSpecial(-4122)
$ws1.Range("A19").Value = 18
$ws1.Range("B19").NumberFormat = "@"
$ws1.Range("B19").Value = "2024.05.04"
$ws1.Range("C19").Value = "广州·黑塔利亚Only"
$ws1.Range("D19").Value = "迎宾大道123号 赛仑吉地大酒店"
$ws1.Range("E19").Value = "2024.05.04 09:30-05.04 16:00"
$ws1.Range("F19").Value = 7
$ws1.Range("G19").Value = 68
$ws1.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=82056"
$ws1.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202402/KI6tnMd81708917202487.jpeg"
# Fix the sequential index column (A) for the rows shifted down by the insert
$ws1.Range("A20").Value = 19
$ws1.Range("A21").Value = 20
$ws1.Range("A22").Value = 21
# Update F20 (shifted old-row19, 萌物语) want-to-go count
$ws1.Range("F20").Value = 516
# F21 (shifted old-row20, 潮宠展) unchanged
# Update F22 (shifted old-row21, 恋与深空only) want-to-go count
$ws1.Range("F22").Value = 465

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 66
$ws2.Range("F6").Value = 15
$ws2.Range("F8").Value = 171
$ws2.Range("F9").Value = 193
$ws2.Range("F10").Value = 42

# --- Sheet 4: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 66
$ws4.Range("F4").Value = 418
$ws4.Range("F7").Value = 573
$ws4.Range("F8").Value = 878
$ws4.Range("F9").Value = 614
$ws4.Range("F10").Value = 790
$ws4.Range("F11").Value = 360
$ws4.Range("F13").Value = 109
$ws4.Range("F14").Value = 1121
$ws4.Range("F15").Value = 570
$ws4.Range("F17").Value = 15
$ws4.Range("F18").Value = 333
$ws4.Range("F19").Value = 448
$ws4.Range("F22").Value = 171
$ws4.Range("F23").Value = 289
$ws4.Range("F26").Value = 193
$ws4.Range("F27").Value = 42

# Insert new row 31 (黑塔利亚Only) in sheet4, shifting old rows 31-33 down to 32-34
$ws4.Rows("31:31").Insert()
$ws4.Range("A32").Copy()
$ws4.Range("A31").PasteSpecial(-4122)
$ws4.Range("B32").Copy()
$ws4.Range("B31").PasteSpecial(-4122)
$ws4.Range("C32:I32").Copy()
$ws4.Range("C31").PasteSpecial(-4122)
$ws4.Range("A31").Value = 30
$ws4.Range("B31").NumberFormat = "@"
$ws4.Range("B31").Value = "2024.05.04"
$ws4.Range("C31").Value = "广州·黑塔利亚Only"
$ws4.Range("D31").Value = "迎宾大道123号 赛仑吉地大酒店"
$ws4.Range("E31").Value = "2024.05.04 09:30-05.04 16:00"
$ws4.Range("F31").Value = 7
$ws4.Range("G31").Value = 68
$ws4.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=82056"
$ws4.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202402/KI6tnMd81708917202487.jpeg"
# Fix the sequential index column (A) for the rows shifted down by the insert
$ws4.Range("A32").Value = 31
$ws4.Range("A33").Value = 32
$ws4.Range("A34").Value = 33
# Update F32 (shifted old-row31, 萌物语) want-to-go count
$ws4.Range("F32").Value = 516
# F33 (shifted old-row32, 潮宠展) unchanged
# Update F34 (shifted old-row33, 恋与深空only) want-to-go count
$ws4.Range("F34").Value = 465
